# daily auto push: 2026-02-02 14:06 UTC
#
# A new timestamp row was inserted into the "sei3" log sheet right after the
# existing 2026/02/02 16:00 entry (row 744), shifting every following row
# (old 744..785) down by one (new 745..786). The new row carries:
#   date=2026/02/02, weekday=月, time=19, ranking=25
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 744; EntireRow.Insert() defaults to shifting the
# existing rows 744+ down (xlShiftDown), exactly like the diff (old row 744
# "2026/12/29" ends up at 745, ..., old row 785 ends up at 786).
$ws.Rows.Item(744).EntireRow.Insert()

# Columns A (date) and B (weekday) in this sheet are plain text ("2026/02/02",
# "月"), not real Excel dates. Force the cells to Text format before writing
# so the literal string is kept verbatim instead of being auto-coerced into a
# date serial number, then drop back to the sheet's normal (General) style so
# no stray per-cell number format is left behind.
$ws.Range("A744:B744").NumberFormat = "@"
$ws.Range("A744").Value = "2026/02/02"
$ws.Range("B744").Value = "月"
$ws.Range("C744").Value = 19
$ws.Range("D744").Value = 25
$ws.Range("A744:B744").Style = "Normal"
